# Generate Report for Handoff
#
# A new handoff run produced a new content GUID / hash pair, so every
# generated file name (and the "latest handoff" timestamps that go with it)
# needs to be refreshed across the Overview, zh-cn and de-de report sheets.
# The hyperlink targets (addresses) themselves are unchanged by this run -
# only the file names shown in the cells (and the hyperlinks' display text)
# and the handoff timestamps change.
#
# NOTE: this engine's Hyperlinks.Delete() clears every hyperlink on the
# whole worksheet (not just the calling range), so each sheet's hyperlinks
# are rebuilt from scratch, in their original order, after the delete.

$wb = $excel.ActiveWorkbook

$oldGuid = "de944284-9911-4ebc-a582-a4766eb4fc4a"
$newGuid = "f9666578-b24b-4dcd-bfd6-c179b4de6882"
$oldHash = "5e720d5b3bf1e98f5d69d71ed4c854d4411140d7"
$newHash = "cd1277526729929c04be5b8889f7d9da818dd730"

$newMdName  = "$newGuid.md"
$newZhCnXlf = "$newGuid.$newHash.zh-cn.xlf"
$newDeDeXlf = "$newGuid.$newHash.de-de.xlf"

# Original hyperlink target addresses (unchanged by this edit).
$addrMd       = "https://github.com/OpenLocalizationTest/oltest/blob/324646384e47dde33b3c7c6e28f4f1e882cd093c/e2e/$oldGuid.md"
$addrZhCnXlf  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e7964d731de29e261ac7bd89ed13cf6f9f537b44/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldGuid.$oldHash.zh-cn.xlf"
$addrDeDeXlf  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d7186aa49f07a55a233730546eed4c44ace134b4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldGuid.$oldHash.de-de.xlf"

# ---------------------------------------------------------------------
# Overview sheet: A2 is a hyperlink whose display text is the .md name;
# D2 is the plain "Latest Handoff Date" text value.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $addrMd, "", "", $newMdName)

$wsOverview.Range("D2").Value = "2016-48-19 08:48:58"

# ---------------------------------------------------------------------
# zh-cn sheet: A2 (.md) and D2 (zh-cn.xlf) are hyperlinks whose display
# text carries the generated file name; E2 is the handoff datetime.
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $addrMd, "", "", $newMdName)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B2"), $addrMd, "", "", ".md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D2"), $addrZhCnXlf, "", "", $newZhCnXlf)

$wsZhCn.Range("E2").Value = "2016-03-19 08:48:55"

# ---------------------------------------------------------------------
# de-de sheet: A2 (.md) and D2 (de-de.xlf) are hyperlinks whose display
# text carries the generated file name; E2 is the handoff datetime.
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $addrMd, "", "", $newMdName)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B2"), $addrMd, "", "", ".md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D2"), $addrDeDeXlf, "", "", $newDeDeXlf)

$wsDeDe.Range("E2").Value = "2016-03-19 08:48:58"
